$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.0167300380228137
$ws.Cells.Item(2, 3).Value = 0.00532319391634981
$ws.Cells.Item(2, 4).Value = 0.00456273764258555
$ws.Cells.Item(2, 5).Value = 0.0372623574144487
$ws.Cells.Item(2, 6).Value = 0.149809885931559
$ws.Cells.Item(2, 7).Value = 0.746768060836502
$ws.Cells.Item(2, 8).Value = 0.0091254752851711
$ws.Cells.Item(2, 9).Value = 0.0212927756653992
$ws.Cells.Item(2, 10).Value = 0.0661596958174905
$ws.Cells.Item(2, 11).Value = 0.96958174904943
$ws.Cells.Item(2, 12).Value = 0.00836501901140684
$ws.Cells.Item(2, 13).Value = 0.00532319391634981
$ws.Cells.Item(2, 14).Value = 0.00608365019011407
$ws.Cells.Item(2, 15).Value = 0
$ws.Cells.Item(2, 16).Value = 0.0174904942965779
$ws.Cells.Item(2, 17).Value = 0.00836501901140684
$ws.Cells.Item(2, 18).Value = 0.109505703422053
$ws.Cells.Item(2, 19).Value = 0.0479087452471483
$ws.Cells.Item(2, 20).Value = 0.0129277566539924
$ws.Cells.Item(2, 21).Value = 0.00456273764258555
$ws.Cells.Item(2, 22).Value = 0.0091254752851711
$ws.Cells.Item(2, 23).Value = 0.0167300380228137
$ws.Cells.Item(2, 24).Value = 0.000760456273764259

$ws.Cells.Item(3, 2).Value = 0.00304182509505703
$ws.Cells.Item(3, 3).Value = 0.0456273764258555
$ws.Cells.Item(3, 4).Value = 0.187072243346008
$ws.Cells.Item(3, 5).Value = 0.942965779467681
$ws.Cells.Item(3, 6).Value = 0.793916349809886
$ws.Cells.Item(3, 7).Value = 0.00684410646387833
$ws.Cells.Item(3, 8).Value = 0.00228136882129278
$ws.Cells.Item(3, 9).Value = 0.00304182509505703
$ws.Cells.Item(3, 10).Value = 0
$ws.Cells.Item(3, 11).Value = 0
$ws.Cells.Item(3, 12).Value = 0.00684410646387833
$ws.Cells.Item(3, 13).Value = 0.000760456273764259
$ws.Cells.Item(3, 14).Value = 0.695817490494297
$ws.Cells.Item(3, 15).Value = 0.00380228136882129
$ws.Cells.Item(3, 16).Value = 0.0365019011406844
$ws.Cells.Item(3, 17).Value = 0
$ws.Cells.Item(3, 18).Value = 0.0182509505703422
$ws.Cells.Item(3, 19).Value = 0.920912547528517
$ws.Cells.Item(3, 20).Value = 0.982509505703422
$ws.Cells.Item(3, 21).Value = 0.00304182509505703
$ws.Cells.Item(3, 22).Value = 0.0152091254752852
$ws.Cells.Item(3, 23).Value = 0.00152091254752852
$ws.Cells.Item(3, 24).Value = 0.00152091254752852

$ws.Cells.Item(4, 2).Value = 0.00304182509505703
$ws.Cells.Item(4, 3).Value = 0.00380228136882129
$ws.Cells.Item(4, 4).Value = 0.00608365019011407
$ws.Cells.Item(4, 5).Value = 0.0159695817490494
$ws.Cells.Item(4, 6).Value = 0.0532319391634981
$ws.Cells.Item(4, 7).Value = 0.244106463878327
$ws.Cells.Item(4, 8).Value = 0.978707224334601
$ws.Cells.Item(4, 9).Value = 0.974904942965779
$ws.Cells.Item(4, 10).Value = 0.926996197718631
$ws.Cells.Item(4, 11).Value = 0.0091254752851711
$ws.Cells.Item(4, 12).Value = 0.984790874524715
$ws.Cells.Item(4, 13).Value = 0.993916349809886
$ws.Cells.Item(4, 14).Value = 0.000760456273764259
$ws.Cells.Item(4, 15).Value = 0.996197718631179
$ws.Cells.Item(4, 16).Value = 0.00836501901140684
$ws.Cells.Item(4, 17).Value = 0.990874524714829
$ws.Cells.Item(4, 18).Value = 0.864638783269962
$ws.Cells.Item(4, 19).Value = 0.000760456273764259
$ws.Cells.Item(4, 20).Value = 0
$ws.Cells.Item(4, 21).Value = 0.0197718631178707
$ws.Cells.Item(4, 22).Value = 0.974144486692015
$ws.Cells.Item(4, 23).Value = 0.979467680608365
$ws.Cells.Item(4, 24).Value = 0.995437262357414

$ws.Cells.Item(5, 2).Value = 0.977186311787072
$ws.Cells.Item(5, 3).Value = 0.945247148288973
$ws.Cells.Item(5, 4).Value = 0.802281368821293
$ws.Cells.Item(5, 5).Value = 0.00380228136882129
$ws.Cells.Item(5, 6).Value = 0.00304182509505703
$ws.Cells.Item(5, 7).Value = 0.00152091254752852
$ws.Cells.Item(5, 8).Value = 0.00988593155893536
$ws.Cells.Item(5, 9).Value = 0.000760456273764259
$ws.Cells.Item(5, 10).Value = 0.00684410646387833
$ws.Cells.Item(5, 11).Value = 0.0212927756653992
$ws.Cells.Item(5, 12).Value = 0
$ws.Cells.Item(5, 13).Value = 0
$ws.Cells.Item(5, 14).Value = 0.297338403041825
$ws.Cells.Item(5, 15).Value = 0
$ws.Cells.Item(5, 16).Value = 0.937642585551331
$ws.Cells.Item(5, 17).Value = 0.000760456273764259
$ws.Cells.Item(5, 18).Value = 0.00760456273764259
$ws.Cells.Item(5, 19).Value = 0.0304182509505703
$ws.Cells.Item(5, 20).Value = 0.00456273764258555
$ws.Cells.Item(5, 21).Value = 0.972623574144487
$ws.Cells.Item(5, 22).Value = 0.000760456273764259
$ws.Cells.Item(5, 23).Value = 0.00228136882129278
$ws.Cells.Item(5, 24).Value = 0

Write-Output "Updated frequency table values for 27341108_gRNA-8"